$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# Insert a new row at row 13, shifting existing rows 13+ down by one.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row with the "classify" key/value pair.
$ws.Range("A13").Value = "classify"
$ws.Range("B13").Value = "CLASSIFY"

# Update the selection to reflect the new cursor position (A14) as in the diff.
$ws.Range("A14").Select()
